$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell J1, matching the style used by the other header cells (A1:I1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "Number total of files"

# Row 2 (Activiti)
$ws.Range("D2").Value = 60
$ws.Range("F2").Value = 49
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = "'    4899`n"
$ws.Range("J2").Style = $ws.Range("C2").Style

# Row 3 (che)
$ws.Range("C3").Value = 39
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 7
$ws.Range("I3").Value = 105
$ws.Range("J3").Value = "'    2538`n"
$ws.Range("J3").Style = $ws.Range("C3").Style

# Row 4 (pinpoint)
$ws.Range("D4").Value = 39
$ws.Range("F4").Value = 43
$ws.Range("J4").Value = "'    8149`n"
$ws.Range("J4").Style = $ws.Range("C4").Style

# Row 5 (skywalking)
$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = "'    2587`n"
$ws.Range("J5").Style = $ws.Range("C5").Style

# Row 6 (wildfly)
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 3
$ws.Range("F6").Value = 170
$ws.Range("G6").Value = 231
$ws.Range("I6").Value = 538
$ws.Range("J6").Value = "'   14079`n"
$ws.Range("J6").Style = $ws.Range("C6").Style

# Row 7 (storm)
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 6
$ws.Range("J7").Value = "'    3398`n"
$ws.Range("J7").Style = $ws.Range("C7").Style
